$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

function Set-PlainValue($cellRef, $val) {
    $ws.Range($cellRef).Value = $val
}

Set-TextValue "D2" "63.783.13"
Set-TextValue "E2" "  +1.13%  "
Set-TextValue "D3" "3.106.97"
Set-TextValue "E3" "  -0.09%  "
Set-TextValue "E4" "  -0.02%  "
Set-TextValue "D5" "584.61"
Set-TextValue "E5" "  -0.01%  "
Set-TextValue "D6" "145.33"
Set-TextValue "E6" "  +0.74%  "
Set-TextValue "E7" "  +0.01%  "
Set-TextValue "D8" "3.103.15"
Set-TextValue "E8" "  +0.03%  "
Set-TextValue "D9" "0.529"
Set-TextValue "E9" "  -0.08%  "
Set-TextValue "E10" "  +7.39%  "
Set-TextValue "E11" "  -2.14%  "
Set-TextValue "D12" "0.458"
Set-TextValue "E12" "  -1.90%  "
Set-TextValue "E13" "  -0.04%  "
Set-TextValue "D14" "36.85"
Set-TextValue "E14" "  +3.75%  "
Set-TextValue "E15" "  -1.26%  "
Set-TextValue "D16" "3.625.55"
Set-TextValue "E16" "  +0.05%  "
Set-TextValue "D17" "63.562.16"
Set-TextValue "E17" "  +0.92%  "
Set-PlainValue "B18" "WrappedEther"
Set-PlainValue "C18" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D18" "3.100.32"
Set-TextValue "E18" "  -0.24%  "
Set-PlainValue "B19" "Polkadot"
Set-PlainValue "C19" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D19" "7.07"
Set-TextValue "E19" "  -1.15%  "
Set-TextValue "D20" "462.26"
Set-TextValue "E20" "  -0.57%  "
Set-TextValue "D21" "14.22"
Set-TextValue "E21" "  +0.29%  "
Set-TextValue "D22" "0.724"
Set-TextValue "E22" "  -0.28%  "
Set-TextValue "D23" "7.44"
Set-TextValue "E23" "  -1.05%  "
Set-TextValue "D24" "12.97"
Set-TextValue "E24" "  -2.30%  "
Set-TextValue "D25" "81.23"
Set-TextValue "E25" "  -0.83%  "
Set-TextValue "E27" "  +0.05%  "
Set-TextValue "D28" "9.20"
Set-TextValue "E28" "  +9.79%  "
Set-TextValue "E29" "  -0.37%  "
Set-TextValue "D30" "2.68"
Set-TextValue "E30" "  +0.17%  "
Set-TextValue "E31" "  -0.04%  "
Set-TextValue "D32" "6.91"
Set-TextValue "E32" "  +0.84%  "
Set-PlainValue "B33" "EthereumClassic"
Set-PlainValue "C33" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D33" "26.73"
Set-TextValue "E33" "  -0.66%  "
Set-PlainValue "B34" "Hedera"
Set-PlainValue "C34" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D34" "0.109"
Set-TextValue "E34" "  -0.37%  "
Set-TextValue "D35" "0.0₃0860"
Set-TextValue "E35" "  -0.24%  "
Set-PlainValue "B36" "dogwifhat"
Set-PlainValue "C36" "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D36" "3.41"
Set-TextValue "E36" "  +2.88%  "
Set-PlainValue "B37" "Stacks"
Set-PlainValue "C37" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D37" "2.32"
Set-TextValue "E37" "  -4.17%  "
Set-PlainValue "B38" "Mantle"
Set-PlainValue "C38" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D38" "1.03"
Set-TextValue "E38" "  -0.48%  "
Set-TextValue "D39" "6.00"
Set-TextValue "E39" "  -0.66%  "
Set-TextValue "D40" "50.39"
Set-TextValue "E40" "  -0.82%  "
Set-TextValue "D41" "437.39"
Set-TextValue "E41" "  +1.19%  "
Set-TextValue "D42" "8.68"
Set-TextValue "E42" "  -0.42%  "
Set-TextValue "D43" "0.0370"
Set-TextValue "E43" "  +0.25%  "
Set-TextValue "D44" "2.880.02"
Set-TextValue "E44" "  -1.48%  "
Set-TextValue "D45" "0.275"
Set-TextValue "E45" "  -1.42%  "
Set-TextValue "E46" "  -2.61%  "
Set-TextValue "D47" "36.31"
Set-TextValue "E47" "  +2.67%  "
Set-TextValue "D48" "125.83"
Set-TextValue "E48" "  +2.36%  "
Set-TextValue "E50" "  -0.61%  "
Set-TextValue "D51" "24.17"
Set-TextValue "E51" "  -1.32%  "
